$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.521.58'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.10%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.287.17'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.43%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '581.49'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.49%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '174.48'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -7.09%  '

$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.579'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.15%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.277.85'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.54%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.173'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -5.53%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.570'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.89%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '45.08'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -4.96%  '

$ws.Range('E13').Value = '  -2.76%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '664.25'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.88%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.825.88'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.16%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.29'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.70%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.556.41'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.14%  '

$ws.Range('E18').Value = '  -0.63%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.296.94'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.21%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.31'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.87%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.79'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.60%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.880'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.19%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.40'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.94%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '16.93'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -6.19%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '98.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.32%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.84'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -4.20%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.64'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -7.20%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.11'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -6.26%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '32.62'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.46%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.32'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -4.46%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.97'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.13%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '573.88'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -6.46%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '10.88'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.09%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.744.55'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.27%  '

$ws.Range('E35').Value = '  -3.68%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.998'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.35'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -13.16%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '55.53'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.87%  '

$ws.Range('E39').Value = '  -1.98%  '

$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '32.13'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.72%  '

$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.60'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -8.68%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.03'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -7.70%  '

$ws.Range('E43').Value = '  -6.55%  '

$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.23'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -5.44%  '

$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.325'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.47%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0400'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -5.49%  '

$ws.Range('E47').Value = '  -0.48%  '

$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.01%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.126'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.45%  '

$ws.Range('E50').Value = '  -1.11%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.77'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.25%  '
